$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)
$ws.Range("C32").Copy()          # style 42 source (unique, no dup earlier since 8 dup is earlier... wait 8<42)
$ws.Range("H5").PasteSpecial(-4122)
